$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3450
$ws.Range("I94").Value = 3266.6667
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 3266.6667
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -2815.6667
$ws.Range("N94").Value = -4902
$ws.Range("H135").Value = 1057.125
$ws.Range("I135").Value = 994.26666
$ws.Range("K135").Value = 8948.399939999999
$ws.Range("M135").Value = -6413.399939999999
$ws.Range("H137").Value = 7145485.5
$ws.Range("I137").Value = 10002150
$ws.Range("J137").Value = 3824.75
$ws.Range("K137").Value = 30006450
$ws.Range("L137").Value = 11474.25
$ws.Range("M137").Value = -30003900
$ws.Range("N137").Value = -16574.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6239.241
$ws.Range("I32").Value = 3593.6143
$ws.Range("J32").Value = 17133
$ws.Range("K32").Value = 3593.6143
$ws.Range("L32").Value = 17133
$ws.Range("M32").Value = -3306.6143
$ws.Range("N32").Value = -17707
$ws.Range("H61").Value = 2876.182
$ws.Range("I61").Value = 1466
$ws.Range("K61").Value = 1466
$ws.Range("M61").Value = -1254
$ws.Range("H88").Value = 1953.4286
$ws.Range("I88").Value = 1953.4286
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1953.4286
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -1547.4286
$ws.Range("H91").Value = 1953.4286
$ws.Range("I91").Value = 1953.4286
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1953.4286
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -549.4286
$ws.Range("H136").Value = 2876.182
$ws.Range("I136").Value = 1466
$ws.Range("K136").Value = 4398
$ws.Range("M136").Value = -1848

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 26813.334
$ws.Range("I86").Value = 2170
$ws.Range("J86").Value = 76100
$ws.Range("K86").Value = 2170
$ws.Range("L86").Value = 76100
$ws.Range("M86").Value = -1047
$ws.Range("N86").Value = -78346
$ws.Range("H89").Value = 26813.334
$ws.Range("I89").Value = 2170
$ws.Range("J89").Value = 76100
$ws.Range("K89").Value = 10850
$ws.Range("L89").Value = 380500
$ws.Range("M89").Value = -5234
$ws.Range("N89").Value = -391732
$ws.Range("H112").Value = 47777.777
$ws.Range("J112").Value = 47777.777
$ws.Range("L112").Value = 47777.777
$ws.Range("N112").Value = -50731.777
$ws.Range("H134").Value = 2289.577
$ws.Range("I134").Value = 1596.619
$ws.Range("K134").Value = 4789.857
$ws.Range("M134").Value = -2254.857

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2043249.8
$ws.Range("I31").Value = 2274391.8
$ws.Range("J31").Value = 9200
$ws.Range("K31").Value = 2274391.8
$ws.Range("L31").Value = 9200
$ws.Range("M31").Value = -2274096.8
$ws.Range("N31").Value = -9790
$ws.Range("H34").Value = 2043249.8
$ws.Range("I34").Value = 2274391.8
$ws.Range("J34").Value = 9200
$ws.Range("K34").Value = 2274391.8
$ws.Range("L34").Value = 9200
$ws.Range("M34").Value = -2274189.8
$ws.Range("N34").Value = -9604
$ws.Range("H132").Value = 2858.742
$ws.Range("I132").Value = 2010.4783
$ws.Range("J132").Value = 5297.5
$ws.Range("K132").Value = 6031.4349
$ws.Range("L132").Value = 15892.5
$ws.Range("M132").Value = -3501.4349
$ws.Range("N132").Value = -20952.5
$ws.Range("H134").Value = 1801.2174
$ws.Range("I134").Value = 1034.5
$ws.Range("K134").Value = 3103.5
$ws.Range("M134").Value = -568.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 190
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 240
$ws.Range("L23").Value = 900
$ws.Range("M23").Value = -5
$ws.Range("N23").Value = -1370
$ws.Range("H68").Value = 3993.889
$ws.Range("J68").Value = 8051
$ws.Range("L68").Value = 24153
$ws.Range("N68").Value = -25775
$ws.Range("H71").Value = 3993.889
$ws.Range("J71").Value = 8051
$ws.Range("L71").Value = 72459
$ws.Range("N71").Value = -80571
$ws.Range("H134").Value = 2840.5
$ws.Range("I134").Value = 1443.375
$ws.Range("J134").Value = 3958.2
$ws.Range("K134").Value = 4330.125
$ws.Range("L134").Value = 11874.6
$ws.Range("M134").Value = 739.875
$ws.Range("N134").Value = -22014.6
$ws.Range("H139").Value = 6328.793
$ws.Range("I139").Value = 2656.1765
$ws.Range("J139").Value = 11531.667
$ws.Range("K139").Value = 7968.529500000001
$ws.Range("L139").Value = 34595.001
$ws.Range("M139").Value = -2828.529500000001
$ws.Range("N139").Value = -44875.001
$ws.Range("H140").Value = 2890.2942
$ws.Range("I140").Value = 960.7143
$ws.Range("J140").Value = 4241
$ws.Range("K140").Value = 2882.1429
$ws.Range("L140").Value = 12723
$ws.Range("M140").Value = 2297.8571
$ws.Range("N140").Value = -23083
$ws.Range("H141").Value = 5186.913
$ws.Range("I141").Value = 3433.1667
$ws.Range("J141").Value = 5805.8823
$ws.Range("K141").Value = 10299.5001
$ws.Range("L141").Value = 17417.6469
$ws.Range("M141").Value = -5119.500100000001
$ws.Range("N141").Value = -27777.6469

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4417
$ws.Range("I70").Value = 4286.6665
$ws.Range("J70").Value = 4612.5
$ws.Range("K70").Value = 4286.6665
$ws.Range("L70").Value = 4612.5
$ws.Range("M70").Value = -4016.6665
$ws.Range("N70").Value = -5152.5
$ws.Range("H73").Value = 4417
$ws.Range("I73").Value = 4286.6665
$ws.Range("J73").Value = 4612.5
$ws.Range("K73").Value = 4286.6665
$ws.Range("L73").Value = 4612.5
$ws.Range("M73").Value = -3350.6665
$ws.Range("N73").Value = -6484.5
$ws.Range("H107").Value = 849.2222
$ws.Range("I107").Value = 341.05554
$ws.Range("J107").Value = 1865.5555
$ws.Range("K107").Value = 341.05554
$ws.Range("L107").Value = 1865.5555
$ws.Range("M107").Value = 1578.94446
$ws.Range("N107").Value = -5705.5555

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 250002750
$ws.Range("I22").Value = 333333660
$ws.Range("K22").Value = 333333660
$ws.Range("M22").Value = -333333365
$ws.Range("H27").Value = 250002750
$ws.Range("I27").Value = 333333660
$ws.Range("K27").Value = 333333660
$ws.Range("M27").Value = -333333553
$ws.Range("H100").Value = 2155
$ws.Range("I100").Value = 1320
$ws.Range("J100").Value = 2751.4285
$ws.Range("K100").Value = 1320
$ws.Range("L100").Value = 2751.4285
$ws.Range("M100").Value = -779
$ws.Range("N100").Value = -3833.4285

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 478281.94
$ws.Range("I122").Value = 716044.9
$ws.Range("K122").Value = 2148134.7
$ws.Range("M122").Value = -2145684.7
$ws.Range("H132").Value = 373208.72
$ws.Range("I132").Value = 557535.4399999999
$ws.Range("J132").Value = 4555.222
$ws.Range("K132").Value = 1672606.32
$ws.Range("L132").Value = 13665.666
$ws.Range("M132").Value = -1670076.32
$ws.Range("N132").Value = -18725.666
